{"js": "// Replace the multiplication expressions in the practice table with the\n// updated problem set. Each old expression is unique in the document, so a\n// targeted search+replace per pair is safe and unambiguous.\nconst replacements = [\n  [\"916\u00d73=\", \"238\u00d79=\"],\n  [\"892\u00d72=\", \"885\u00d76=\"],\n  [\"741\u00d72=\", \"589\u00d76=\"],\n  [\"830\u00d79=\", \"734\u00d73=\"],\n  [\"294\u00d73=\", \"685\u00d78=\"],\n  [\"333\u00d78=\", \"592\u00d79=\"],\n  [\"547\u00d72=\", \"476\u00d79=\"],\n  [\"291\u00d77=\", \"923\u00d77=\"],\n  [\"858\u00d72=\", \"843\u00d73=\"],\n  [\"855\u00d73=\", \"827\u00d74=\"],\n  [\"485\u00d74=\", \"486\u00d72=\"],\n  [\"318\u00d77=\", \"785\u00d72=\"],\n  [\"781\u00d76=\", \"699\u00d72=\"],\n  [\"224\u00d72=\", \"736\u00d73=\"],\n  [\"364\u00d75=\", \"723\u00d78=\"],\n  [\"140\u00d79=\", \"354\u00d77=\"],\n  [\"464\u00d74=\", \"524\u00d74=\"],\n  [\"308\u00d76=\", \"761\u00d76=\"],\n  [\"546\u00d79=\", \"407\u00d73=\"],\n  [\"348\u00d78=\", \"366\u00d79=\"],\n  [\"111\u00d77=\", \"125\u00d78=\"],\n  [\"611\u00d75=\", \"825\u00d79=\"],\n  [\"128\u00d75=\", \"146\u00d74=\"],\n  [\"382\u00d79=\", \"330\u00d72=\"],\n  [\"601\u00d78=\", \"152\u00d72=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"text\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the multiplication-practice table: each old expression is unique\n# in the document, so Find/Replace per pair is unambiguous.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"916\u00d73=\", \"238\u00d79=\"),\n    @(\"892\u00d72=\", \"885\u00d76=\"),\n    @(\"741\u00d72=\", \"589\u00d76=\"),\n    @(\"830\u00d79=\", \"734\u00d73=\"),\n    @(\"294\u00d73=\", \"685\u00d78=\"),\n    @(\"333\u00d78=\", \"592\u00d79=\"),\n    @(\"547\u00d72=\", \"476\u00d79=\"),\n    @(\"291\u00d77=\", \"923\u00d77=\"),\n    @(\"858\u00d72=\", \"843\u00d73=\"),\n    @(\"855\u00d73=\", \"827\u00d74=\"),\n    @(\"485\u00d74=\", \"486\u00d72=\"),\n    @(\"318\u00d77=\", \"785\u00d72=\"),\n    @(\"781\u00d76=\", \"699\u00d72=\"),\n    @(\"224\u00d72=\", \"736\u00d73=\"),\n    @(\"364\u00d75=\", \"723\u00d78=\"),\n    @(\"140\u00d79=\", \"354\u00d77=\"),\n    @(\"464\u00d74=\", \"524\u00d74=\"),\n    @(\"308\u00d76=\", \"761\u00d76=\"),\n    @(\"546\u00d79=\", \"407\u00d73=\"),\n    @(\"348\u00d78=\", \"366\u00d79=\"),\n    @(\"111\u00d77=\", \"125\u00d78=\"),\n    @(\"611\u00d75=\", \"825\u00d79=\"),\n    @(\"128\u00d75=\", \"146\u00d74=\"),\n    @(\"382\u00d79=\", \"330\u00d72=\"),\n    @(\"601\u00d78=\", \"152\u00d72=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Replacement.ClearFormatting()\n    $rng.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
